# TC46_Canine_Filter_Breed-WestHlnd.xlsx -- "corrected ICDC Breed 1-14 scripts"
#
# The FilesTab query (row 4, column B) is corrected: the `File Type` and
# `Breed` projection lines are dropped from the Cypher/RETURN clause, and the
# row is re-sized to match its new (shorter) line count. The sheet view is
# also scrolled/selected down to row 4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newFilesQuery = @'
MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (f)-[*]->(c:case)<--(demo:demographic)
WHERE demo.breed IN ['West Highland White Terrier'] 
OPTIONAL MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
WITH DISTINCT f, parent, c, demo, diag, s
RETURN  coalesce(f.file_name, '') AS `File Name`,
        coalesce(labels(parent)[0], '') AS `Association`,
        coalesce(f.file_description, '') AS `Description`,
        coalesce(f.file_format, '') AS `Format`,
        coalesce(f.file_size, '') AS `Size`,
        coalesce(c.case_id, '') AS `Case ID`,
         coalesce(diag.disease_term,'') AS Diagnosis , 
        coalesce(s.clinical_study_designation,'') AS `Study Code`
'@

$ws.Range("B4").Value = $newFilesQuery

# The corrected query has two fewer lines, so the wrapped row shrinks.
$ws.Rows.Item(4).RowHeight = 217.5

# Scroll/select down to the row that was just edited.
[void]$ws.Range("B4").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 4
$win.ScrollColumn = 1
